# Mole Rat stat-block update
# - "Burrower" trait renamed to "Burrowing" and rewritten with new burrow rules
#   (also splits its description into several runs).
# - "Natural Weapons" trait's description runs are collapsed back into one run.
# - The Description paragraph gets a spell-check proofErr wrapper around
#   "molerat" (splitting that single run into three).

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Edit 1: "Burrower." -> "Burrowing." trait rewrite
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Burrower. The mole rat has a burrowing speed of 15 ft. It costs the mole rat 5 feet of movement to begin burrowing within a pre-existing entry point.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$xml1 = "<w:p xmlns:w='$wNs'>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space='preserve'>Burrowing. </w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>The </w:t></w:r>" +
        "<w:r><w:t>mole rat</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> has a burrowing speed of 10 feet through loose earth and 0 feet through solid rock and metal.</w:t></w:r>" +
        "</w:p>"
$rng.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: "Natural Weapons." trait - merge the three plain runs into one
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Natural Weapons. The mole rat" + [char]8217 + "s unarmed attacks use a d6 instead of a d4.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$xml2 = "<w:p xmlns:w='$wNs'>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space='preserve'>Natural Weapons. </w:t></w:r>" +
        "<w:r><w:t>The mole rat" + [char]8217 + "s unarmed attacks use a d6 instead of a d4.</w:t></w:r>" +
        "</w:p>"
$rng.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Edit 3: Description paragraph - wrap "molerat" with spell-check proofErr
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Mole rats are mutated rodents much larger than their pre-War counterpart, the naked molerat. They have survived as a species by burrowing underground where the soil protected them from direct nuclear explosions. However, they still were greatly mutated by the subsequent fallout, increasing their size as well as their viciousness. A few people in the wastes have even been known to train them and keep them as pets.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$xml3 = "<w:p xmlns:w='$wNs'>" +
        "<w:r><w:t xml:space='preserve'>Mole rats are mutated rodents much larger than their pre-War counterpart, the naked </w:t></w:r>" +
        "<w:proofErr w:type='spellStart'/>" +
        "<w:r><w:t>molerat</w:t></w:r>" +
        "<w:proofErr w:type='spellEnd'/>" +
        "<w:r><w:t>. They have survived as a species by burrowing underground where the soil protected them from direct nuclear explosions. However, they still were greatly mutated by the subsequent fallout, increasing their size as well as their viciousness. A few people in the wastes have even been known to train them and keep them as pets.</w:t></w:r>" +
        "</w:p>"
$rng.InsertXML($xml3)
